$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B2:B118 accuracy values per updated epoch run
$newValues = @(0.921875,0.875,0.875,0.84375,0.8125,0.78125,0.734375,0.71875,0.71875,0.65625,0.6875,0.703125,0.6875,0.640625,0.625,0.765625,0.671875,0.625,0.578125,0.5625,0.46875,0.53125,0.515625,0.53125,0.53125,0.53125,0.53125,0.53125,0.53125,0.53125,0.53125,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.546875,0.53125,0.53125,0.53125,0.53125,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.515625,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.625,0.6875,0.484375,0.59375,0.5625,0.640625,0.609375,0.5625,0.578125,0.59375,0.609375,0.640625,0.59375,0.59375,0.625,0.459016393442623)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $newValues[$i]
}

# Updated object id text for the DisplayOutputs repr in column A (rows 102-118)
$newAddr = "<__main__.DisplayOutputs object at 0x7f10cc11ce20>"
for ($r = 102; $r -le 118; $r++) {
    $ws.Cells.Item($r, 1).Value = $newAddr
}

$ws.Range("A2:B118").Select()
